$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing angle_diff (D) values for rows 4-7, and SOA (C) values too
$ws.Range("C4").Value = -2
$ws.Range("C5").Value = -2
$ws.Range("C6").Value = -1
$ws.Range("D6").Value = -1
$ws.Range("C7").Value = -1
$ws.Range("D7").Value = -1

# Append new data rows (8-13) for the stim2_c = 0 condition block
$newData = @(
    @(0,  0, -0.4, -0.4),
    @(12, 0, -0.4, -0.4),
    @(0,  0, -2,   -2),
    @(12, 0, -2,   -2),
    @(0,  0, -1,   -1),
    @(12, 0, -1,   -1)
)

$startRow = 8
for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $newData[$i][0]
    $ws.Range("B$row").Value = $newData[$i][1]
    $ws.Range("C$row").Value = $newData[$i][2]
    $ws.Range("D$row").Value = $newData[$i][3]
}

$ws.Range("B8:B13").Select()
